$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.245.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.718.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06181"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.717.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07072"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5970"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.422"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.258.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006796"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.937.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.539"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.726"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.286"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.400"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.759"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.971"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.676"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07759"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04460"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.95%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6180"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9252"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "113.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +17.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.423"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.917"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01483"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.568"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3818"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1181"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.277"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.788"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3373"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  +1.41%  "
